# Release MHD 4.2.2 close #419
#
# Updates the "Metadata" sheet of the MHDlistTypes CodeSystem workbook:
#   - Version bumped 4.2.1 -> 4.2.2
#   - Date updated to the new publication date
#   - Contact block expanded from a single "No display for ContactDetail"
#     placeholder into the three rendered ContactDetail lines

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "4.2.2"

# Date
$ws.Range("B8").Value = "2024-05-18T12:39:23-05:00"

# Contact (rows 10-12 already exist as three "Contact" rows; replace the
# placeholder text in column B with the resolved ContactDetail lines)
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
